$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the coin price/volume figures for this run.
# We restore each cell's original Style afterwards so that only the text
# content changes (Excel would otherwise reinterpret digit/dot strings such
# as "215.14" as numbers and silently round/alter their display).

$orig_D2 = $ws.Range("D2").Style
$orig_E2 = $ws.Range("E2").Style
$orig_D3 = $ws.Range("D3").Style
$orig_E3 = $ws.Range("E3").Style
$orig_E4 = $ws.Range("E4").Style
$orig_D5 = $ws.Range("D5").Style
$orig_E5 = $ws.Range("E5").Style
$orig_E6 = $ws.Range("E6").Style
$orig_E7 = $ws.Range("E7").Style
$orig_D8 = $ws.Range("D8").Style
$orig_E8 = $ws.Range("E8").Style
$orig_D9 = $ws.Range("D9").Style
$orig_E9 = $ws.Range("E9").Style
$orig_E10 = $ws.Range("E10").Style
$orig_D11 = $ws.Range("D11").Style
$orig_E11 = $ws.Range("E11").Style
$orig_D12 = $ws.Range("D12").Style
$orig_E12 = $ws.Range("E12").Style
$orig_D13 = $ws.Range("D13").Style
$orig_E14 = $ws.Range("E14").Style
$orig_E15 = $ws.Range("E15").Style
$orig_D16 = $ws.Range("D16").Style
$orig_E16 = $ws.Range("E16").Style
$orig_D17 = $ws.Range("D17").Style
$orig_E17 = $ws.Range("E17").Style
$orig_D18 = $ws.Range("D18").Style
$orig_E18 = $ws.Range("E18").Style
$orig_D19 = $ws.Range("D19").Style
$orig_E19 = $ws.Range("E19").Style
$orig_D20 = $ws.Range("D20").Style
$orig_E20 = $ws.Range("E20").Style
$orig_E21 = $ws.Range("E21").Style
$orig_E22 = $ws.Range("E22").Style
$orig_D23 = $ws.Range("D23").Style
$orig_E23 = $ws.Range("E23").Style
$orig_E24 = $ws.Range("E24").Style
$orig_D25 = $ws.Range("D25").Style
$orig_E25 = $ws.Range("E25").Style
$orig_D26 = $ws.Range("D26").Style
$orig_E26 = $ws.Range("E26").Style
$orig_D27 = $ws.Range("D27").Style
$orig_E27 = $ws.Range("E27").Style
$orig_E28 = $ws.Range("E28").Style
$orig_E29 = $ws.Range("E29").Style
$orig_E30 = $ws.Range("E30").Style
$orig_D31 = $ws.Range("D31").Style
$orig_E31 = $ws.Range("E31").Style
$orig_D32 = $ws.Range("D32").Style
$orig_E32 = $ws.Range("E32").Style
$orig_D33 = $ws.Range("D33").Style
$orig_E33 = $ws.Range("E33").Style
$orig_E34 = $ws.Range("E34").Style
$orig_E35 = $ws.Range("E35").Style
$orig_D36 = $ws.Range("D36").Style
$orig_E36 = $ws.Range("E36").Style
$orig_E37 = $ws.Range("E37").Style
$orig_D38 = $ws.Range("D38").Style
$orig_E38 = $ws.Range("E38").Style
$orig_E39 = $ws.Range("E39").Style
$orig_D40 = $ws.Range("D40").Style
$orig_E40 = $ws.Range("E40").Style
$orig_E42 = $ws.Range("E42").Style
$orig_D43 = $ws.Range("D43").Style
$orig_E43 = $ws.Range("E43").Style
$orig_E44 = $ws.Range("E44").Style
$orig_D45 = $ws.Range("D45").Style
$orig_E45 = $ws.Range("E45").Style
$orig_D46 = $ws.Range("D46").Style
$orig_E46 = $ws.Range("E46").Style
$orig_D47 = $ws.Range("D47").Style
$orig_E47 = $ws.Range("E47").Style
$orig_D48 = $ws.Range("D48").Style
$orig_E48 = $ws.Range("E48").Style
$orig_E49 = $ws.Range("E49").Style
$orig_E50 = $ws.Range("E50").Style
$orig_D51 = $ws.Range("D51").Style
$orig_E51 = $ws.Range("E51").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.638.94'
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").Value = '1.668.65'
$ws.Range("E3").Value = '  -3.16%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '215.14'
$ws.Range("E5").Value = '  -1.73%  '
$ws.Range("E6").Value = '  -2.29%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '24.03'
$ws.Range("E8").Value = '  -2.05%  '
$ws.Range("D9").Value = '0.263'
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("E10").Value = '  -1.72%  '
$ws.Range("D11").Value = '0.0879'
$ws.Range("E11").Value = '  -2.00%  '
$ws.Range("D12").Value = '1.904.22'
$ws.Range("E12").Value = '  -3.18%  '
$ws.Range("D13").Value = '1.667.81'
$ws.Range("E14").Value = '  -3.05%  '
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").Value = '66.74'
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("D17").Value = '27.613.91'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").Value = '243.58'
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").Value = '0.0₃0731'
$ws.Range("E19").Value = '  -3.21%  '
$ws.Range("D20").Value = '7.70'
$ws.Range("E20").Value = '  -4.09%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("E22").Value = '  -2.91%  '
$ws.Range("D23").Value = '9.34'
$ws.Range("E23").Value = '  -3.90%  '
$ws.Range("E24").Value = '  -4.33%  '
$ws.Range("D25").Value = '147.35'
$ws.Range("E25").Value = '  -1.15%  '
$ws.Range("D26").Value = '7.22'
$ws.Range("E26").Value = '  -3.61%  '
$ws.Range("D27").Value = '16.53'
$ws.Range("E27").Value = '  -1.14%  '
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("E29").Value = '  -2.28%  '
$ws.Range("E30").Value = '  +2.86%  '
$ws.Range("D31").Value = '0.0503'
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("D32").Value = '3.36'
$ws.Range("E32").Value = '  -2.40%  '
$ws.Range("D33").Value = '1.471.49'
$ws.Range("E33").Value = '  -1.64%  '
$ws.Range("E34").Value = '  -4.71%  '
$ws.Range("E35").Value = '  -4.96%  '
$ws.Range("D36").Value = '2.37'
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("E37").Value = '  -2.80%  '
$ws.Range("D38").Value = '0.576'
$ws.Range("E38").Value = '  -5.08%  '
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("D40").Value = '69.66'
$ws.Range("E40").Value = '  -1.56%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").Value = '5.42'
$ws.Range("E43").Value = '  -7.34%  '
$ws.Range("E44").Value = '  -2.77%  '
$ws.Range("D45").Value = '1.811.93'
$ws.Range("E45").Value = '  -3.12%  '
$ws.Range("D46").Value = '0.787'
$ws.Range("E46").Value = '  -1.88%  '
$ws.Range("D47").Value = '1.75'
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("D48").Value = '89.30'
$ws.Range("E48").Value = '  -1.85%  '
$ws.Range("E49").Value = '  -4.17%  '
$ws.Range("E50").Value = '  -2.02%  '
$ws.Range("D51").Value = '7.94'
$ws.Range("E51").Value = '  -3.83%  '

$ws.Range("D2").Style = $orig_D2
$ws.Range("E2").Style = $orig_E2
$ws.Range("D3").Style = $orig_D3
$ws.Range("E3").Style = $orig_E3
$ws.Range("E4").Style = $orig_E4
$ws.Range("D5").Style = $orig_D5
$ws.Range("E5").Style = $orig_E5
$ws.Range("E6").Style = $orig_E6
$ws.Range("E7").Style = $orig_E7
$ws.Range("D8").Style = $orig_D8
$ws.Range("E8").Style = $orig_E8
$ws.Range("D9").Style = $orig_D9
$ws.Range("E9").Style = $orig_E9
$ws.Range("E10").Style = $orig_E10
$ws.Range("D11").Style = $orig_D11
$ws.Range("E11").Style = $orig_E11
$ws.Range("D12").Style = $orig_D12
$ws.Range("E12").Style = $orig_E12
$ws.Range("D13").Style = $orig_D13
$ws.Range("E14").Style = $orig_E14
$ws.Range("E15").Style = $orig_E15
$ws.Range("D16").Style = $orig_D16
$ws.Range("E16").Style = $orig_E16
$ws.Range("D17").Style = $orig_D17
$ws.Range("E17").Style = $orig_E17
$ws.Range("D18").Style = $orig_D18
$ws.Range("E18").Style = $orig_E18
$ws.Range("D19").Style = $orig_D19
$ws.Range("E19").Style = $orig_E19
$ws.Range("D20").Style = $orig_D20
$ws.Range("E20").Style = $orig_E20
$ws.Range("E21").Style = $orig_E21
$ws.Range("E22").Style = $orig_E22
$ws.Range("D23").Style = $orig_D23
$ws.Range("E23").Style = $orig_E23
$ws.Range("E24").Style = $orig_E24
$ws.Range("D25").Style = $orig_D25
$ws.Range("E25").Style = $orig_E25
$ws.Range("D26").Style = $orig_D26
$ws.Range("E26").Style = $orig_E26
$ws.Range("D27").Style = $orig_D27
$ws.Range("E27").Style = $orig_E27
$ws.Range("E28").Style = $orig_E28
$ws.Range("E29").Style = $orig_E29
$ws.Range("E30").Style = $orig_E30
$ws.Range("D31").Style = $orig_D31
$ws.Range("E31").Style = $orig_E31
$ws.Range("D32").Style = $orig_D32
$ws.Range("E32").Style = $orig_E32
$ws.Range("D33").Style = $orig_D33
$ws.Range("E33").Style = $orig_E33
$ws.Range("E34").Style = $orig_E34
$ws.Range("E35").Style = $orig_E35
$ws.Range("D36").Style = $orig_D36
$ws.Range("E36").Style = $orig_E36
$ws.Range("E37").Style = $orig_E37
$ws.Range("D38").Style = $orig_D38
$ws.Range("E38").Style = $orig_E38
$ws.Range("E39").Style = $orig_E39
$ws.Range("D40").Style = $orig_D40
$ws.Range("E40").Style = $orig_E40
$ws.Range("E42").Style = $orig_E42
$ws.Range("D43").Style = $orig_D43
$ws.Range("E43").Style = $orig_E43
$ws.Range("E44").Style = $orig_E44
$ws.Range("D45").Style = $orig_D45
$ws.Range("E45").Style = $orig_E45
$ws.Range("D46").Style = $orig_D46
$ws.Range("E46").Style = $orig_E46
$ws.Range("D47").Style = $orig_D47
$ws.Range("E47").Style = $orig_E47
$ws.Range("D48").Style = $orig_D48
$ws.Range("E48").Style = $orig_E48
$ws.Range("E49").Style = $orig_E49
$ws.Range("E50").Style = $orig_E50
$ws.Range("D51").Style = $orig_D51
$ws.Range("E51").Style = $orig_E51
